$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# style 3: FFBE94A3 -> FF99B27B
$ws.Range("I2:O2").Interior.Color = 8106649

# style 4: FF81FC8C -> FFC3C2AE
$ws.Range("C3:N3").Interior.Color = 11453123

# style 5: FF819E7A -> FFF482A4
$ws.Range("C4:J4").Interior.Color = 10781428

# style 6: FFF787D8 -> FFDAEF80
$ws.Range("C5:O5").Interior.Color = 8450010
$ws.Range("C10:K10").Interior.Color = 8450010
$ws.Range("L13:N13").Interior.Color = 8450010
$ws.Range("C14:K14").Interior.Color = 8450010
$ws.Range("L18:N18").Interior.Color = 8450010
$ws.Range("C19:K19").Interior.Color = 8450010
$ws.Range("L34:O34").Interior.Color = 8450010
$ws.Range("C35:K35").Interior.Color = 8450010

# style 7: FFF6DBFB -> FFD9D3F8
$ws.Range("C6:E6").Interior.Color = 16307161

# style 8: FFB481DC -> FF91FE81
$ws.Range("C7:E7").Interior.Color = 8519313

# style 9: FFD5C2E7 -> FF7CF6C6
$ws.Range("C8:K8").Interior.Color = 13039228
$ws.Range("C18:D18").Interior.Color = 13039228
$ws.Range("C31:K31").Interior.Color = 13039228
$ws.Range("C34:K34").Interior.Color = 13039228

# style 10: FFD9CCAD -> FF7ED194
$ws.Range("L8:N8").Interior.Color = 9752958
$ws.Range("C9:K9").Interior.Color = 9752958
$ws.Range("C24:D24").Interior.Color = 9752958

# style 11: FF8AB5FD -> FFC9FAD7
$ws.Range("C11:K11").Interior.Color = 14154441

# style 12: FFFC7B96 -> FFFBE4B3
$ws.Range("C12:D12").Interior.Color = 11789563
$ws.Range("C30:D30").Interior.Color = 11789563

# style 13: FFBBD07A -> FFB3DA8B
$ws.Range("C13:K13").Interior.Color = 9165491

# style 14: FF94E9B9 -> FFB0B2DD
$ws.Range("C15:J15").Interior.Color = 14529200
$ws.Range("K16").Interior.Color = 14529200

# style 15: FFE6A183 -> FFD787E2
$ws.Range("K15").Interior.Color = 14845911
$ws.Range("C16:J16").Interior.Color = 14845911
$ws.Range("C26:E26").Interior.Color = 14845911
$ws.Range("L26:O26").Interior.Color = 14845911
$ws.Range("F27:K27").Interior.Color = 14845911

# style 16: FF9FCBE0 -> FFFA7EFC
$ws.Range("C17").Interior.Color = 16547578

# style 17: FF7EF3F8 -> FF849DAC
$ws.Range("E18:K18").Interior.Color = 11312516
$ws.Range("C20:O20").Interior.Color = 11312516

# style 18: FF8AA6BF -> FFF7FDEC
$ws.Range("C21:D21").Interior.Color = 15531511
$ws.Range("C22:D22").Interior.Color = 15531511

# style 19: FF937F99 -> FF9ED8FD
$ws.Range("C23:D23").Interior.Color = 16636062

# style 20: FFD1F9C9 -> FF7F80EE
$ws.Range("C25:E25").Interior.Color = 15630463
$ws.Range("L25:O25").Interior.Color = 15630463
$ws.Range("F26:K26").Interior.Color = 15630463
$ws.Range("C28:K28").Interior.Color = 15630463

# style 21: FFF5EE7D -> FFB57C83
$ws.Range("F25:K25").Interior.Color = 8617141

# style 22: FFD19AFE -> FFEBB284
$ws.Range("C29:K29").Interior.Color = 8696555

# style 23: FFA2B995 -> FFA37BC8
$ws.Range("C32:F32").Interior.Color = 13138851
$ws.Range("C33:I33").Interior.Color = 13138851

# style 24: FF7B80F4 -> FFFEC1EB
$ws.Range("C36:J36").Interior.Color = 15450622

